$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "LED" row (A13:C13, merged B13:C13) needs to move up to row 10 while
# a fresh, still-tall-but-empty row 13 is left behind (waiting for the
# upcoming led-show module content).

# Unmerge first so the cut/paste below carries plain, individually
# addressable cells.
$ws.Range("B13:C13").UnMerge()

# Move the whole block (values + formatting) from row 13 to row 10.
$ws.Range("A13:C13").Cut($ws.Range("A10:C10"))

# Re-create the merge at the new location.
$ws.Range("B10:C10").Merge()

# Give the relocated row its own (shorter) height.
$ws.Rows.Item(10).RowHeight = 51.75

# Wipe whatever is left behind at row 13 ...
$ws.Range("A13:C13").Clear()
# ... but keep the row reserved with its original tall height.
$ws.Rows.Item(13).RowHeight = 117

# Match the author's final cursor position.
$ws.Range("E11").Select()
